$wb = $excel.ActiveWorkbook

# --- Duplicate "Nädal 2" into a new "Nädal 3" sheet placed right after it ---
$ws2 = $wb.Worksheets.Item("Nädal 2")
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Nädal 3"

# --- Clear out the copied week's log rows (keep the row-number column A) ---
$ws3.Range("B8:J18").ClearContents()

# --- Log a new entry in row 7 (12.02.2020, 22:30-23:40, 70 min, GitHub) ---
$ws3.Range("B7").Value = 43873
$ws3.Range("C7").Value = 0.9375
$ws3.Range("D7").Value = 0.98611111111111116
$ws3.Range("E7").ClearContents()
$ws3.Range("F7").Value = 70
$ws3.Range("G7").Value = "GitHub"
$ws3.Range("H7").Value = "üritan commiti tagasi saada"
$ws3.Range("I7").ClearContents()
$ws3.Range("J7").ClearContents()

# --- Selections: Nädal 2 no longer active, whole table selected; Nädal 3 is now active, F8 selected ---
[void]$ws2.Range("A1:J19").Select()
[void]$ws3.Range("F8").Select()
[void]$ws3.Activate()
